$wb = $excel.ActiveWorkbook

# Work on the HIV sheet (second sheet), which is the active/tabbed sheet
$ws = $wb.Worksheets.Item("HIV")

# Add the missing discrepancy text into D5, matching the style/content
# already used in D2:D4 (quote-prefixed text starting with "-")
$ws.Range("D5").Value = "'-HIV Drawn Date is more than 2 yrs.:"

# Update the active selection on the HIV sheet
$ws.Range("D7").Select()
